$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23..75 down to 24..76.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new weekly record.
$ws.Cells.Item(23, 1).Value = 4
$ws.Cells.Item(23, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(23, 3).Value = "Los Lagos"
$ws.Cells.Item(23, 4).Value = 45133
$ws.Cells.Item(23, 5).Value = 10
$ws.Cells.Item(23, 6).Value = 100112012
$ws.Cells.Item(23, 7).Value = "Espinaca"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 25
$ws.Cells.Item(23, 11).Value = 14000
$ws.Cells.Item(23, 12).Value = 14000
$ws.Cells.Item(23, 13).Value = 14000
$ws.Cells.Item(23, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(23, 15).Value = "Región Metropolitana"
$ws.Cells.Item(23, 16).Value = 1400
$ws.Cells.Item(23, 17).Value = 10
$ws.Cells.Item(23, 18).Value = "Hortaliza"
